$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 386.58
$ws.Cells.Item(92, 9).Value = 295.92105
$ws.Cells.Item(92, 10).Value = 673.6667
$ws.Cells.Item(92, 11).Value = 295.92105
$ws.Cells.Item(92, 12).Value = 673.6667
$ws.Cells.Item(92, 13).Value = 952.0789500000001
$ws.Cells.Item(92, 14).Value = -3169.6667
$ws.Cells.Item(133, 8).Value = 38800
$ws.Cells.Item(133, 10).Value = 38800
$ws.Cells.Item(133, 12).Value = 38800
$ws.Cells.Item(133, 14).Value = -48920
$ws.Cells.Item(135, 8).Value = 879.7436
$ws.Cells.Item(135, 9).Value = 541.5833
$ws.Cells.Item(135, 10).Value = 4937.6665
$ws.Cells.Item(135, 11).Value = 4874.2497
$ws.Cells.Item(135, 12).Value = 44438.9985
$ws.Cells.Item(135, 13).Value = -2339.2497
$ws.Cells.Item(135, 14).Value = -49508.9985
$ws.Cells.Item(137, 8).Value = 1921.0952
$ws.Cells.Item(137, 9).Value = 1439.1052
$ws.Cells.Item(137, 10).Value = 6500
$ws.Cells.Item(137, 11).Value = 4317.3156
$ws.Cells.Item(137, 12).Value = 19500
$ws.Cells.Item(137, 13).Value = -1767.3156
$ws.Cells.Item(137, 14).Value = -24600
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 63437.688
$ws.Cells.Item(2, 9).Value = 845.8182
$ws.Cells.Item(2, 10).Value = 201139.8
$ws.Cells.Item(2, 11).Value = 845.8182
$ws.Cells.Item(2, 12).Value = 201139.8
$ws.Cells.Item(2, 13).Value = -732.8182
$ws.Cells.Item(2, 14).Value = -201365.8
$ws.Cells.Item(4, 8).Value = 175
$ws.Cells.Item(4, 9).Value = 250
$ws.Cells.Item(4, 11).Value = 250
$ws.Cells.Item(4, 13).Value = -134
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 14).Value = ""
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).Value = ""
$ws.Cells.Item(23, 8).Value = 56803.8
$ws.Cells.Item(32, 8).Value = 32042.322
$ws.Cells.Item(32, 9).Value = 5094.5107
$ws.Cells.Item(32, 11).Value = 5094.5107
$ws.Cells.Item(32, 13).Value = -4807.5107
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = ""
$ws.Cells.Item(37, 14).Value = ""
$ws.Cells.Item(44, 8).Value = 12857.25
$ws.Cells.Item(44, 10).Value = 12857.25
$ws.Cells.Item(44, 12).Value = 12857.25
$ws.Cells.Item(44, 14).Value = -13833.25
$ws.Cells.Item(55, 8).Value = 11017.667
$ws.Cells.Item(55, 10).Value = 13026.5
$ws.Cells.Item(55, 12).Value = 13026.5
$ws.Cells.Item(55, 14).Value = -13656.5
$ws.Cells.Item(61, 8).Value = 1779.1666
$ws.Cells.Item(61, 9).Value = 1212.0834
$ws.Cells.Item(61, 10).Value = 2913.3333
$ws.Cells.Item(61, 11).Value = 1212.0834
$ws.Cells.Item(61, 12).Value = 2913.3333
$ws.Cells.Item(61, 13).Value = -1000.0834
$ws.Cells.Item(61, 14).Value = -3337.3333
$ws.Cells.Item(80, 8).Value = 26697.2
$ws.Cells.Item(80, 10).Value = 26697.2
$ws.Cells.Item(80, 12).Value = 26697.2
$ws.Cells.Item(80, 14).Value = -28693.2
$ws.Cells.Item(83, 8).Value = 26697.2
$ws.Cells.Item(83, 10).Value = 26697.2
$ws.Cells.Item(83, 12).Value = 80091.60000000001
$ws.Cells.Item(83, 14).Value = -90075.60000000001
$ws.Cells.Item(116, 8).Value = 63437.688
$ws.Cells.Item(116, 9).Value = 845.8182
$ws.Cells.Item(116, 10).Value = 201139.8
$ws.Cells.Item(116, 11).Value = 845.8182
$ws.Cells.Item(116, 12).Value = 201139.8
$ws.Cells.Item(116, 13).Value = 1448.1818
$ws.Cells.Item(116, 14).Value = -205727.8
$ws.Cells.Item(136, 8).Value = 1779.1666
$ws.Cells.Item(136, 9).Value = 1212.0834
$ws.Cells.Item(136, 10).Value = 2913.3333
$ws.Cells.Item(136, 11).Value = 3636.2502
$ws.Cells.Item(136, 12).Value = 8739.999899999999
$ws.Cells.Item(136, 13).Value = -1086.2502
$ws.Cells.Item(136, 14).Value = -13839.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 63437.688
$ws.Cells.Item(3, 9).Value = 845.8182
$ws.Cells.Item(3, 10).Value = 201139.8
$ws.Cells.Item(3, 11).Value = 845.8182
$ws.Cells.Item(3, 12).Value = 201139.8
$ws.Cells.Item(3, 13).Value = -731.8182
$ws.Cells.Item(3, 14).Value = -201367.8
$ws.Cells.Item(86, 8).Value = 45037.42
$ws.Cells.Item(86, 9).Value = 102127.27
$ws.Cells.Item(86, 10).Value = 3171.5334
$ws.Cells.Item(86, 11).Value = 102127.27
$ws.Cells.Item(86, 12).Value = 3171.5334
$ws.Cells.Item(86, 13).Value = -101004.27
$ws.Cells.Item(86, 14).Value = -5417.5334
$ws.Cells.Item(89, 8).Value = 45037.42
$ws.Cells.Item(89, 9).Value = 102127.27
$ws.Cells.Item(89, 10).Value = 3171.5334
$ws.Cells.Item(89, 11).Value = 510636.35
$ws.Cells.Item(89, 12).Value = 15857.667
$ws.Cells.Item(89, 13).Value = -505020.35
$ws.Cells.Item(89, 14).Value = -27089.667
$ws.Cells.Item(134, 8).Value = 2350.982
$ws.Cells.Item(134, 9).Value = 2121.1064
$ws.Cells.Item(134, 11).Value = 6363.3192
$ws.Cells.Item(134, 13).Value = -3828.3192
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 21091.143
$ws.Cells.Item(19, 9).Value = 46
$ws.Cells.Item(19, 10).Value = 36875
$ws.Cells.Item(19, 11).Value = 46
$ws.Cells.Item(19, 12).Value = 36875
$ws.Cells.Item(19, 13).Value = 124
$ws.Cells.Item(19, 14).Value = -37215
$ws.Cells.Item(24, 8).Value = 21091.143
$ws.Cells.Item(24, 9).Value = 46
$ws.Cells.Item(24, 10).Value = 36875
$ws.Cells.Item(24, 11).Value = 46
$ws.Cells.Item(24, 12).Value = 36875
$ws.Cells.Item(24, 13).Value = 124
$ws.Cells.Item(24, 14).Value = -37215
$ws.Cells.Item(31, 8).Value = 31564.963
$ws.Cells.Item(31, 9).Value = 1006.8
$ws.Cells.Item(31, 10).Value = 50085.062
$ws.Cells.Item(31, 11).Value = 1006.8
$ws.Cells.Item(31, 12).Value = 50085.062
$ws.Cells.Item(31, 13).Value = -711.8
$ws.Cells.Item(31, 14).Value = -50675.062
$ws.Cells.Item(34, 8).Value = 31564.963
$ws.Cells.Item(34, 9).Value = 1006.8
$ws.Cells.Item(34, 10).Value = 50085.062
$ws.Cells.Item(34, 11).Value = 1006.8
$ws.Cells.Item(34, 12).Value = 50085.062
$ws.Cells.Item(34, 13).Value = -804.8
$ws.Cells.Item(34, 14).Value = -50489.062
$ws.Cells.Item(132, 8).Value = 27780788
$ws.Cells.Item(132, 9).Value = 25002744
$ws.Cells.Item(132, 10).Value = 35718056
$ws.Cells.Item(132, 11).Value = 75008232
$ws.Cells.Item(132, 12).Value = 107154168
$ws.Cells.Item(132, 13).Value = -75005702
$ws.Cells.Item(132, 14).Value = -107159228
$ws.Cells.Item(133, 8).Value = 68000
$ws.Cells.Item(133, 10).Value = 68000
$ws.Cells.Item(133, 12).Value = 68000
$ws.Cells.Item(133, 14).Value = -73060
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1473.2433
$ws.Cells.Item(5, 9).Value = 671.3333
$ws.Cells.Item(5, 10).Value = 2953.6924
$ws.Cells.Item(5, 11).Value = 2013.9999
$ws.Cells.Item(5, 12).Value = 8861.0772
$ws.Cells.Item(5, 13).Value = -1901.9999
$ws.Cells.Item(5, 14).Value = -9085.0772
$ws.Cells.Item(17, 8).Value = 286.66666
$ws.Cells.Item(17, 9).Value = 250
$ws.Cells.Item(17, 10).Value = 305
$ws.Cells.Item(17, 11).Value = 750
$ws.Cells.Item(17, 12).Value = 915
$ws.Cells.Item(17, 13).Value = -581
$ws.Cells.Item(17, 14).Value = -1253
$ws.Cells.Item(34, 8).Value = 651.4286
$ws.Cells.Item(34, 10).Value = 1050
$ws.Cells.Item(34, 12).Value = 3150
$ws.Cells.Item(34, 14).Value = -3318
$ws.Cells.Item(39, 8).Value = 3100
$ws.Cells.Item(39, 10).Value = 3100
$ws.Cells.Item(39, 12).Value = 9300
$ws.Cells.Item(39, 14).Value = -9888
$ws.Cells.Item(55, 8).Value = 14816.25
$ws.Cells.Item(55, 10).Value = 10450.77
$ws.Cells.Item(55, 12).Value = 31352.31
$ws.Cells.Item(55, 14).Value = -31706.31
$ws.Cells.Item(58, 8).Value = 2699.6
$ws.Cells.Item(58, 9).Value = 1999.5
$ws.Cells.Item(58, 10).Value = 3166.3333
$ws.Cells.Item(58, 11).Value = 5998.5
$ws.Cells.Item(58, 12).Value = 9498.999899999999
$ws.Cells.Item(58, 13).Value = -5870.5
$ws.Cells.Item(58, 14).Value = -9754.999899999999
$ws.Cells.Item(122, 8).Value = 444
$ws.Cells.Item(122, 10).Value = 397
$ws.Cells.Item(122, 12).Value = 3573
$ws.Cells.Item(122, 14).Value = -8473
$ws.Cells.Item(131, 8).Value = 838.14
$ws.Cells.Item(131, 10).Value = 877.0323
$ws.Cells.Item(131, 12).Value = 2631.0969
$ws.Cells.Item(131, 14).Value = -12711.0969
$ws.Cells.Item(135, 8).Value = 1473.2433
$ws.Cells.Item(135, 9).Value = 671.3333
$ws.Cells.Item(135, 10).Value = 2953.6924
$ws.Cells.Item(135, 11).Value = 6041.9997
$ws.Cells.Item(135, 12).Value = 26583.2316
$ws.Cells.Item(135, 13).Value = -3506.9997
$ws.Cells.Item(135, 14).Value = -31653.2316
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 152
$ws.Cells.Item(2, 9).Value = 154.09091
$ws.Cells.Item(2, 10).Value = 148.16667
$ws.Cells.Item(2, 11).Value = 154.09091
$ws.Cells.Item(2, 12).Value = 148.16667
$ws.Cells.Item(2, 13).Value = -41.09091000000001
$ws.Cells.Item(2, 14).Value = -374.16667
$ws.Cells.Item(70, 8).Value = 99452.95
$ws.Cells.Item(70, 9).Value = 157401.08
$ws.Cells.Item(70, 11).Value = 157401.08
$ws.Cells.Item(70, 13).Value = -157131.08
$ws.Cells.Item(73, 8).Value = 99452.95
$ws.Cells.Item(73, 9).Value = 157401.08
$ws.Cells.Item(73, 11).Value = 157401.08
$ws.Cells.Item(73, 13).Value = -156465.08
$ws.Cells.Item(113, 8).Value = 1288.8
$ws.Cells.Item(113, 9).Value = 737
$ws.Cells.Item(113, 10).Value = 1656.6666
$ws.Cells.Item(113, 11).Value = 737
$ws.Cells.Item(113, 12).Value = 1656.6666
$ws.Cells.Item(113, 13).Value = 1433
$ws.Cells.Item(113, 14).Value = -5996.6666
$ws.Cells.Item(126, 8).Value = 5886033
$ws.Cells.Item(126, 9).Value = 4500.6665
$ws.Cells.Item(126, 10).Value = 14708332
$ws.Cells.Item(126, 11).Value = 13501.9995
$ws.Cells.Item(126, 12).Value = 44124996
$ws.Cells.Item(126, 13).Value = -11031.9995
$ws.Cells.Item(126, 14).Value = -44129936
$ws.Cells.Item(132, 8).Value = 2804.5
$ws.Cells.Item(132, 9).Value = 1848.1852
$ws.Cells.Item(132, 11).Value = 5544.5556
$ws.Cells.Item(132, 13).Value = -3014.5556
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8466.666999999999
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 8466.666999999999
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 8466.666999999999
$ws.Cells.Item(7, 13).Value = ""
$ws.Cells.Item(7, 14).Value = -8690.666999999999
$ws.Cells.Item(100, 8).Value = 2354.8333
$ws.Cells.Item(100, 9).Value = 1999.6
$ws.Cells.Item(100, 11).Value = 1999.6
$ws.Cells.Item(100, 13).Value = -1458.6
$ws.Cells.Item(122, 8).Value = 2134
$ws.Cells.Item(122, 9).Value = 1400
$ws.Cells.Item(122, 10).Value = 2280.8
$ws.Cells.Item(122, 11).Value = 4200
$ws.Cells.Item(122, 12).Value = 6842.400000000001
$ws.Cells.Item(122, 13).Value = -1750
$ws.Cells.Item(122, 14).Value = -11742.4
$ws.Cells.Item(126, 8).Value = 8466.666999999999
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 8466.666999999999
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 25400.001
$ws.Cells.Item(126, 13).Value = ""
$ws.Cells.Item(126, 14).Value = -30340.001
$ws.Cells.Item(132, 8).Value = 3327.3142
$ws.Cells.Item(132, 9).Value = 3076.8438
$ws.Cells.Item(132, 11).Value = 9230.5314
$ws.Cells.Item(132, 13).Value = -6700.5314
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1992.7241
$ws.Cells.Item(132, 9).Value = 1777.7451
$ws.Cells.Item(132, 10).Value = 3559
$ws.Cells.Item(132, 11).Value = 5333.2353
$ws.Cells.Item(132, 12).Value = 10677
$ws.Cells.Item(132, 13).Value = -2803.2353
$ws.Cells.Item(132, 14).Value = -15737
